$d = $word.ActiveDocument

$found0 = $d.Content.Find.Execute('Kweli nqanaba nikezela ngeefomu zemvume. Tell participants to detach the last page and sign it if they still consent to being part of the study. They can fold the last page in half and hand it back to you as you go around (don’t let them pass it between themselves). Please remind them that all the information they need about the study is on the two pages they still have and that they should keep those safely as it also has contact numbers on it. ', $true, $false, $false, $false, $false, $true, 1, $false, 'Kweli nqanaba nikezela ngeefomu zemvume. Xelela abathathi-nxaxheba ukuba bakhuphe iphepha lokugqibela kwaye balityikitye ukuba basavuma ukuba yinxalenye yophononongo. Basenokulisonga iphepha lokugqibela kwihafu kwaye balibuyisele kuwe njengoko ujikeleza (ungabavumeli ukuba bawadlulise phakathi kwabo). Nceda ubakhumbuze ukuba lonke ulwazi abalifunayo malunga nophononongo likumaphepha amabini abasenawo kwaye kufuneka bawagcine ngokukhuselekileyo njengoko enawo neenombolo zoqhagamshelwano kuwo. ', 2)
if (-not $found0) { Write-Output "NOT FOUND: 0" }

$found1 = $d.Content.Find.Execute('Thanks again for participating in this discussion. Please let me know if any of my questions are confusing.', $true, $false, $false, $false, $false, $true, 1, $false, 'Enkosi kwakhona ngothatha inxaxheba kule ngxoxo. Nceda undazise ukuba nayiphi na imibuzo yam iyakubhida.', 2)
if (-not $found1) { Write-Output "NOT FOUND: 1" }

$found2 = $d.Content.Find.Execute('I’d like to start by asking you how you found out about the ParentText programme.', $true, $false, $false, $false, $false, $true, 1, $false, 'Ndingathanda ukuqala ngokukubuza ukuba ufumanise njani malunga nenkqubo ye-ParentText.', 2)
if (-not $found2) { Write-Output "NOT FOUND: 2" }

$found3 = $d.Content.Find.Execute('How did you find out about or hear about ParentText? (E.g. you saw a poster, someone you know shared it with you, a nurse at the clinic told you about it)', $true, $false, $false, $false, $false, $true, 1, $false, 'Ufumanise njani malunga okanye uve njani nge-ParentText? (Umz. ubone ipowusta, umntu omaziyo wabelana ngayo nawe, umongikazi wasekliniki ukuxelele ngayo)', 2)
if (-not $found3) { Write-Output "NOT FOUND: 3" }

$found4 = $d.Content.Find.Execute('For those who saw the poster: What did you think of the poster?', $true, $false, $false, $false, $false, $true, 1, $false, 'Kwabo bathe babona ipowusta: Ucinge ntoni nge powusta?', 2)
if (-not $found4) { Write-Output "NOT FOUND: 4" }

$found5 = $d.Content.Find.Execute('Probe: Did it catch your attention? Did you get all of the information that you needed from it? Did you like the layout e.g. images, colours?', $true, $false, $false, $false, $false, $true, 1, $false, 'Buza: Ingaba ikutsalile? Ingaba ulufumene lonke ulwazi obulifuna kuyo? Ingaba uluthandile ulwakhiwo umz. imifanekiso, imibala?', 2)
if (-not $found5) { Write-Output "NOT FOUND: 5" }

$found6 = $d.Content.Find.Execute('For those who were told about it by a nurse: What was your experience of that? What was it that made you want to sign up?', $true, $false, $false, $false, $false, $true, 1, $false, 'Kwabo abathe baxelelwa ngumongikazi malunga nayo: Ebenjani amava akho ngoku? What was it that made you want to sign up?', 2)
if (-not $found6) { Write-Output "NOT FOUND: 6" }

$found7 = $d.Content.Find.Execute('Did any of you find out about it in another way? - besides through the clinic?', $true, $false, $false, $false, $false, $true, 1, $false, 'Ingaba abanye benu bayazi malunga nayo ngenye indlela? - ngaphandle kokuyazi nge klinikhi?', 2)
if (-not $found7) { Write-Output "NOT FOUND: 7" }

$found8 = $d.Content.Find.Execute('Is there anything else you’d like to share about how you found out about or heard about the ParentText programme?', $true, $false, $false, $false, $false, $true, 1, $false, 'Ingaba ikhona enye into ongathanda ukwabelana ngayo malunga nendlela ofumanise ngayo okanye ove ngayo malunga nenkqubo ye-ParentText?', 2)
if (-not $found8) { Write-Output "NOT FOUND: 8" }

$found9 = $d.Content.Find.Execute('B. First impression of ParentText: ', $true, $false, $false, $false, $false, $true, 1, $false, 'B. Imbonakalo yokuqala ye-ParentText: ', 2)
if (-not $found9) { Write-Output "NOT FOUND: 9" }

$found10 = $d.Content.Find.Execute('Let''s talk about your first thoughts of the ParentText programme. I know it was 3 months ago or more for some of you since you first started interacting with it. So, I want to just take you back. Once you joined ParentText you would have started with some research questions, and you would have received your first thank you in the form of R25 airtime. Once you finished that you would have started interacting with the actual lessons and chatbot….', $true, $false, $false, $false, $false, $true, 1, $false, 'Masikhe sithethe malunga neengcinga zakho zokuqala ngenkqubo ye-ParentText. Ndiyayazi ukuba ibizinyanga ezintathu ezidlulileyo okanye nangaphezulu kwabanye benu oko waqala ukunxibelelana nayo. Ke, ndifuna ukukubuyisela emva. Nje ukuba ujoyine i-ParentText ubuya kuqala ngemibuzo yophando, kwaye ubuya kufumana umbulelo wakho wokuqala we-R25 airtime. Ukugqiba nje kwakho oko ubuqalile ukunxibelelana ngako nezona zifundo kunye ne-chatbot….', 2)
if (-not $found10) { Write-Output "NOT FOUND: 10" }

$found11 = $d.Content.Find.Execute('What was your first impression of ParentText? ', $true, $false, $false, $false, $false, $true, 1, $false, 'Yayiyintoni umbono wakho wokuqala we-ParentText? ', 2)
if (-not $found11) { Write-Output "NOT FOUND: 11" }

$found12 = $d.Content.Find.Execute('Probe: What are the things you liked? What worked well for you? How was it helpful for you? ', $true, $false, $false, $false, $false, $true, 1, $false, 'Buza: Zintoni izinto ozithandileyo? Yintoni ekusebenzele kakuhle wena? Ibeluncedo njani kuwe? ', 2)
if (-not $found12) { Write-Output "NOT FOUND: 12" }

$found13 = $d.Content.Find.Execute('Probe: What didn’t you like? How can we make it better?', $true, $false, $false, $false, $false, $true, 1, $false, 'Buza: Yintoni ongayithandanga? Singayenza njani ibengcono?', 2)
if (-not $found13) { Write-Output "NOT FOUND: 13" }

$found14 = $d.Content.Find.Execute('What were you expecting when you started using ParentText? ', $true, $false, $false, $false, $false, $true, 1, $false, 'Ubulindele ntoni ukuqala kwakho ukusebenzisa i-ParentText? ', 2)
if (-not $found14) { Write-Output "NOT FOUND: 14" }

$found15 = $d.Content.Find.Execute('Probe: Did it meet your expectations?', $true, $false, $false, $false, $false, $true, 1, $false, 'Buza: Ikufezekisile obukulindele?', 2)
if (-not $found15) { Write-Output "NOT FOUND: 15" }

$found16 = $d.Content.Find.Execute('Probe: How did it meet your expectations?', $true, $false, $false, $false, $false, $true, 1, $false, 'Buza: Ikufezekise njani obukulindele?', 2)
if (-not $found16) { Write-Output "NOT FOUND: 16" }

$found17 = $d.Content.Find.Execute('Probe: What was it that didn’t meet your expectations?', $true, $false, $false, $false, $false, $true, 1, $false, 'Buza: Yintoni engakhange ifikelele kobukulindele wena?', 2)
if (-not $found17) { Write-Output "NOT FOUND: 17" }

$found18 = $d.Content.Find.Execute('3. When first connecting to and using ParentText you would have needed WhatsApp to be connected. What was your experience of connecting to and using ParentText?', $true, $false, $false, $false, $false, $true, 1, $false, '3. Xa uqala ukuqhagamshelana kwaye usebenzisa i-ParentText ubuya kudinga u-WhatsApp ukuze uqhagamshelane. Athini amava akho ngokuqhagamshelana nokusebenzisa i-ParentText?', 2)
if (-not $found18) { Write-Output "NOT FOUND: 18" }

$found19 = $d.Content.Find.Execute('Probe: How were you all connected? E.g. mobile data, Wi-Fi, a combination. ', $true, $false, $false, $false, $false, $true, 1, $false, 'Buza: Beniqhagamshelana njani nonke? Umz. idata yefowuni, Wi-Fi, indibaniselwano. ', 2)
if (-not $found19) { Write-Output "NOT FOUND: 19" }

$found20 = $d.Content.Find.Execute('Probe: What challenges did you have with connecting?', $true, $false, $false, $false, $false, $true, 1, $false, 'Buza: Zeziphi iingxaki oye wadibana nazo ngokuqhagamshelana?', 2)
if (-not $found20) { Write-Output "NOT FOUND: 20" }
